$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '42.525.49'
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").Value = '2.289.99'
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue "D5" '156.70'
$ws.Range("E5").Value = '  +15,555.08%  '
Set-TextValue "D6" '307.24'
$ws.Range("E6").Value = '  +1.19%  '
Set-TextValue "D7" '96.80'
$ws.Range("E7").Value = '  +5.94%  '
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  +3.87%  '
Set-TextValue "D11" '36.33'
$ws.Range("E11").Value = '  +13.33%  '
Set-TextValue "D12" '0.0805'
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("E13").Value = '  -1.49%  '
Set-TextValue "D14" '6.73'
$ws.Range("E14").Value = '  +2.63%  '
$ws.Range("D15").Value = '2.641.39'
$ws.Range("E15").Value = '  +1.33%  '
Set-TextValue "D16" '14.68'
$ws.Range("E16").Value = '  +3.56%  '
$ws.Range("D17").Value = '2.292.53'
$ws.Range("E17").Value = '  +0.21%  '
Set-TextValue "D18" '0.806'
$ws.Range("D19").Value = '42.394.52'
$ws.Range("E19").Value = '  +2.13%  '
Set-TextValue "D20" '12.86'
$ws.Range("E20").Value = '  +4.61%  '
$ws.Range("D21").Value = '0.0₃0920'
$ws.Range("E21").Value = '  +1.91%  '
Set-TextValue "D22" '6.02'
$ws.Range("E22").Value = '  +2.20%  '
Set-TextValue "D23" '67.85'
$ws.Range("E23").Value = '  +1.90%  '
Set-TextValue "D24" '243.21'
$ws.Range("E24").Value = '  +1.31%  '
Set-TextValue "D25" '2.61'
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  +2.66%  '
$ws.Range("E27").Value = '  -0.17%  '
Set-TextValue "D28" '24.03'
$ws.Range("E28").Value = '  +0.42%  '
Set-TextValue "D29" '36.64'
$ws.Range("E29").Value = '  +7.58%  '
Set-TextValue "D30" '9.61'
$ws.Range("E30").Value = '  +1.27%  '
Set-TextValue "D31" '2.10'
$ws.Range("E31").Value = '  +1.98%  '
Set-TextValue "D32" '161.04'
$ws.Range("E32").Value = '  +0.64%  '
Set-TextValue "D33" '5.34'
$ws.Range("E33").Value = '  +3.31%  '
Set-TextValue "D34" '0.999'
$ws.Range("E34").Value = '  -0.06%  '
Set-TextValue "D35" '0.0752'
$ws.Range("E35").Value = '  +1.65%  '
$ws.Range("E36").Value = '  +3.15%  '
Set-TextValue "D37" '17.46'
$ws.Range("E37").Value = '  +5.32%  '
$ws.Range("E38").Value = '  +3.57%  '
$ws.Range("E40").Value = '  +0.50%  '
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("E42").Value = '  +6.22%  '
Set-TextValue "D43" '2.42'
$ws.Range("E43").Value = '  +18.77%  '
$ws.Range("D44").Value = '2.004.74'
$ws.Range("E44").Value = '  -2.12%  '
Set-TextValue "D45" '19.40'
$ws.Range("E45").Value = '  +0.53%  '
Set-TextValue "D46" '0.0286'
$ws.Range("E46").Value = '  +3.16%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D47" '3.03'
$ws.Range("E47").Value = '  +5.89%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D48" '10.26'
$ws.Range("E48").Value = '  -0.58%  '
Set-TextValue "D49" '54.09'
$ws.Range("E49").Value = '  +4.98%  '
Set-TextValue "D50" '1.56'
$ws.Range("E50").Value = '  +1.79%  '
Set-TextValue "D51" '72.74'
$ws.Range("E51").Value = '  +0.26%  '
